# This script updates the "SnippetID" column (column H) values in the
# active worksheet. These were auto-generated 4-character codes that
# need to be regenerated/replaced with new values, row by row, from row 2
# through row 61 (row 1 is the header "SnippetID" and is left unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSnippetIds = @(
    "Uboh", "h2eg", "HLVR", "qCE6", "1Y8A", "Ry2W", "Swb0", "GjFV",
    "Ibtw", "Ibtw", "Ibtw", "Ibtw", "w0Sk", "5Cz6", "YnCs", "7ASK",
    "7ASK", "zuUX", "MK6X", "i0Kg", "FaBu", "6ypN", "qtLT", "sYCw",
    "M4ZR", "33ji", "W7YK", "sn1B", "zPtx", "f8LK", "pINk", "ONfP",
    "9Y1E", "9Y1E", "C2h0", "C2h0", "6E7M", "6E7M", "2HPo", "2HPo",
    "bhDy", "bhDy", "ORi7", "ORi7", "0rPJ", "0rPJ", "1Ysb", "uRBE",
    "uRBE", "xoT7", "xoT7", "HgC0", "HgC0", "HgC0", "w4KQ", "w4KQ",
    "qdWm", "p7nG", "73iL", "7nV5"
)

$startRow = 2
for ($i = 0; $i -lt $newSnippetIds.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("H$row").Value = $newSnippetIds[$i]
}
